$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: was the "abc" bot row, now holds the raw Telegram chat id
$ws.Range("A2").Value = 1082171472
$ws.Range("B2").Value = 2
$ws.Range("H2").Value = 3
$ws.Range("I2").Value = 1

# Row 3: "zyx" placeholder renamed to "test1"
$ws.Range("A3").Value = "test1"
$ws.Range("H3").Value = 1

# Remove the now-unused 4th row entirely (its data merged/removed)
$ws.Rows.Item(4).Delete()
